# UKIM resources update: rewrite the header row of the "people responsible"
# upload template with the new wording/order used for the movement-of-goods
# template, and leave the cursor on A2 (matching the author's resave).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Full name of person responsible for the movement of goods"
$ws.Range("B1").Value = "Residential address"
$ws.Range("C1").Value = "Date of birth"
$ws.Range("D1").Value = "National Insurance number"
$ws.Range("E1").Value = "identification number if no National Insurance number  (eg passport number, driver's licence, national identity card)"

# Column A grew to fit the new, longer header text.
$ws.Columns.Item(1).AutoFit()

# Leave the selection on A2, as in the saved workbook.
$ws.Range("A2").Select()
